$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix team name typos / formatting in the "Equipos" (teams) sheet.
$ws.Range("A4").Value = "Trataguas Team"
$ws.Range("B45:B60").Value = "Liga Veteranos +30"
$ws.Range("A76").Value = "Valleseco Vet. B F7"

# Last user interaction left the cursor on A77 (single cell, not a range).
$ws.Range("A77").Select()
